# Apply timetable refactor: H2/H9 swap + evening/Saturday-morning reassignment.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row block 11:00-12:45 (rows 7-14), columns B-E ---
# Old merge layout: B7:B14, C7:C10, C11:C14, D7:D10, D11:D14, E7:E14
# New merge layout: B7:B10, B11:B14, C7:C14, D7:D14, E7:E10, E11:E14
$ws.Range("B7:B14").UnMerge()
$ws.Range("C7:C10").UnMerge()
$ws.Range("C11:C14").UnMerge()
$ws.Range("D7:D10").UnMerge()
$ws.Range("D11:D14").UnMerge()
$ws.Range("E7:E14").UnMerge()

$ws.Range("B7").Value = "Practice " + [char]10 + "(Harp practice room)"
$ws.Range("B11").Value = "Private Lesson with Sivan MEGAN " + [char]10 + "(Room 245)"
$ws.Range("C7").Value = "Free Time"
$ws.Range("C11").Value = ""
$ws.Range("D7").Value = "Practice " + [char]10 + "(Harp practice room)"
$ws.Range("D11").Value = ""
$ws.Range("E7").Value = "Practice " + [char]10 + "(Harp practice room)"
$ws.Range("E11").Value = "Private Lesson with Sivan MEGAN " + [char]10 + "(Room 245)"

$ws.Range("B7:B10").Merge()
$ws.Range("B11:B14").Merge()
$ws.Range("C7:C14").Merge()
$ws.Range("D7:D14").Merge()
$ws.Range("E7:E10").Merge()
$ws.Range("E11:E14").Merge()

# --- Row 20 (14:15 block) ---
$ws.Range("B20").Value = "Practice " + [char]10 + "(Harp practice room)"
$ws.Range("C20").Value = "Free Time"
$ws.Range("D20").Value = "Private Lesson with Gwyneth WENTINK " + [char]10 + "(Room 236)"
$ws.Range("E20").Value = "Practice " + [char]10 + "(Harp practice room)"

# --- Row 24 (15:15 block) : Acting Class -> Ensemble (Room 245) ---
$ws.Range("B24").Value = "Ensemble " + [char]10 + "(Room 245)"
$ws.Range("C24").Value = "Ensemble " + [char]10 + "(Room 245)"
$ws.Range("D24").Value = "Ensemble " + [char]10 + "(Room 245)"
$ws.Range("E24").Value = "Ensemble " + [char]10 + "(Room 245)"
$ws.Range("F24").Value = "Ensemble " + [char]10 + "(Room 245)"

# --- Row 28 (16:15 block) : Ensemble (Room 236) -> Acting Class / Break -> Group Activity ---
$ws.Range("B28").Value = "Acting Class " + [char]10 + "(Room G13)"
$ws.Range("C28").Value = "Acting Class " + [char]10 + "(Room G13)"
$ws.Range("D28").Value = "Acting Class " + [char]10 + "(Room G13)"
$ws.Range("E28").Value = "Acting Class " + [char]10 + "(Room G13)"
$ws.Range("F28").Value = "Group Activity " + [char]10 + "(Room Group Activity)"
